# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Source cells are plain text ("Price"/"Volume(1h)" columns), so numeric-looking
# price strings are written with a leading quote to force text, then the cell
# style is reset to Normal so no stray number-format/quote-prefix style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.895.97'
$ws.Range("E2").Value = '  -1.86%  '

$ws.Range("D3").Value = '3.325.65'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = "'575.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '

$ws.Range("D6").Value = "'182.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").Value = '3.326.45'
$ws.Range("E9").Value = '  +1.71%  '

$ws.Range("E10").Value = '  -0.71%  '

$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '3.906.46'
$ws.Range("E13").Value = '  +1.72%  '

$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").Value = "'27.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").Value = '67.117.09'
$ws.Range("E16").Value = '  -1.51%  '

$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("D18").Value = '3.323.02'
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").Value = "'442.27"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'13.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.97%  '

$ws.Range("D21").Value = "'5.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.07%  '

$ws.Range("E22").Value = '  +2.24%  '

$ws.Range("D23").Value = "'73.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.92%  '

$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("D25").Value = '3.477.45'
$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("D26").Value = "'0.511"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.62%  '

$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("E28").Value = '  +3.75%  '

$ws.Range("D29").Value = "'8.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.34%  '

$ws.Range("D30").Value = "'1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("D31").Value = "'1.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.30%  '

$ws.Range("D32").Value = "'22.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = "'5.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.42%  '

$ws.Range("E35").Value = '  -0.75%  '

$ws.Range("E36").Value = '  -1.82%  '

$ws.Range("D37").Value = "'162.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.57%  '

$ws.Range("E38").Value = '  +4.10%  '

$ws.Range("D39").Value = "'27.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.63%  '

$ws.Range("E40").Value = '  -2.30%  '

$ws.Range("D41").Value = '2.826.73'
$ws.Range("E41").Value = '  +7.41%  '

$ws.Range("D42").Value = "'0.788"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("E44").Value = '  -1.31%  '

$ws.Range("D45").Value = "'40.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("D47").Value = "'24.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.32%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = "'321.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.07%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = "'2.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.67%  '

$ws.Range("E50").Value = '  -0.31%  '

$ws.Range("D51").Value = "'0.983"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.45%  '
